$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose data (columns B..AB) must be swapped with each other.
# Column A (the sequential record index) stays untouched.
$pairs = @(
    @(9, 10),
    @(29, 30),
    @(49, 50),
    @(76, 77),
    @(111, 112),
    @(122, 123),
    @(177, 178)
)

$firstCol = 2   # column B
$lastCol = 28   # column AB

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $row1Values = @{}
    $row2Values = @{}

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $row1Values[$c] = $ws.Cells.Item($r1, $c).Value2
        $row2Values[$c] = $ws.Cells.Item($r2, $c).Value2
    }

    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r1, $c).Value = $row2Values[$c]
        $ws.Cells.Item($r2, $c).Value = $row1Values[$c]
    }
}
